{"js": "// [Track-344] AS combine amendment sections\n//\n// Before: \"...: {$item_val[i].project_description}\" ends one run/paragraph,\n// and a whole separate paragraph holds \"{$item_val[i].project_description:showEnd}\".\n// After: both placeholders live in the same run/paragraph:\n// \"...: {$item_val[i].project_description}{$item_val[i].project_description:showEnd}\"\n// and the now-redundant paragraph is removed.\n\nconst body = context.document.body;\n\n// The literal \"[i].\" (not \"[i+1].\") pins this to the first amendment block,\n// which is the one the diff touches.\nconst hits = body.search(\"{$item_val[i].project_description}\", {\n  matchCase: true,\n  matchWildcards: false\n});\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find target placeholder '{$item_val[i].project_description}'\");\n}\n\nconst hit = hits.items[0];\nhit.load(\"text,paragraphs\");\nawait context.sync();\n\nconst firstParagraph = hit.paragraphs.getFirst();\nconst secondParagraph = firstParagraph.getNext();\nsecondParagraph.load(\"text\");\nawait context.sync();\n\nconst mergedText = hit.text + secondParagraph.text;\n\n// Replace the placeholder text in-place (keeps it inside the existing run,\n// preserving that run's formatting) so it now also carries the showEnd tag.\nhit.insertText(mergedText, \"Replace\");\n\n// The showEnd text now lives twice: once merged above, and once in the old,\n// now-redundant, paragraph. Delete that paragraph (merges it away).\nsecondParagraph.delete();\n\nawait context.sync();\n", "ps1": "# [Track-344] AS combine amendment sections\n#\n# Before: \"...: {$item_val[i].project_description}\" ends one run/paragraph,\n# and a whole separate paragraph holds \"{$item_val[i].project_description:showEnd}\".\n# After: both placeholders live in the same run/paragraph:\n# \"...: {$item_val[i].project_description}{$item_val[i].project_description:showEnd}\"\n# and the now-redundant paragraph is removed.\n\n$d = $word.ActiveDocument\n\n# The literal \"[i].\" (not \"[i+1].\") pins this to the first amendment block,\n# which is the one the diff touches.\n$target = \"{`$item_val[i].project_description}\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $target\n$found = $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false)\nif (-not $found) {\n    throw \"Could not find target placeholder '$target'\"\n}\n\n# Locate which paragraph (by 1-based index) contains the found range, so we\n# can reach its very-next paragraph (the one holding the ':showEnd' tag).\n$allParagraphs = $d.Paragraphs\n$paragraphIndex = -1\nfor ($i = 1; $i -le $allParagraphs.Count; $i++) {\n    $p = $allParagraphs.Item($i)\n    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {\n        $paragraphIndex = $i\n        break\n    }\n}\nif ($paragraphIndex -eq -1) {\n    throw \"Could not locate the paragraph containing '$target'\"\n}\n\n$firstParagraph = $allParagraphs.Item($paragraphIndex)\n$secondParagraph = $allParagraphs.Item($paragraphIndex + 1)\n\n# A paragraph's Range.Text carries its trailing paragraph mark (chr 13) and,\n# for the last paragraph in a table cell, the cell-end mark (chr 7) too --\n# strip those so only the visible placeholder text is reused.\n$tailText = $secondParagraph.Range.Text.TrimEnd([char]13, [char]7, [char]12)\n\n# Append the \"showEnd\" placeholder text onto the existing run that ends with\n# \"{$item_val[i].project_description}\" (Find/Replace rewrites the matched\n# text in place, so it stays inside that run with its original formatting).\n$find2 = $rng.Find\n$find2.ClearFormatting()\n$find2.Text = $target\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = $target + $tailText\n$find2.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, ($target + $tailText), 2)\n\n# The showEnd text now lives twice: once merged above, and once in the old,\n# now-redundant, second paragraph. Delete that paragraph (merges it away).\n$secondParagraph.Range.Delete()\n"}
